# Insert a new price-record row at position 445 (pushes the existing
# rows 445-544 down to 446-545, extending the used range to A1:R545).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(445).Insert()

$ws.Range("A445").Value = 9
$ws.Range("B445").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C445").Value = 'Metropolitana'
$ws.Range("D445").Value = 44943
$ws.Range("E445").Value = 13
$ws.Range("F445").Value = 100112012
$ws.Range("G445").Value = 'Espinaca'
$ws.Range("H445").Value = 'Sin especificar'
$ws.Range("I445").Value = 'Primera'
$ws.Range("J445").Value = 160
$ws.Range("K445").Value = 6000
$ws.Range("L445").Value = 8000
$ws.Range("M445").Value = 7000
$ws.Range("N445").Value = '$/cuna 10 kilos'
$ws.Range("O445").Value = 'Provincia de Chacabuco'
$ws.Range("P445").Value = 700
$ws.Range("Q445").Value = 10
$ws.Range("R445").Value = 'Hortaliza'
